# Template test update:
#  - Fill B2:K2 with the new "ab[colmeta]" label (adds a shared string).
#  - Widen column B and split the old 3-11 run so column J gets its own
#    (slightly different) width.
#  - Move the active selection to I16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, columns B..K: tag them with the "ab[colmeta]" label.
$ws.Range("B2:K2").Value = "ab[colmeta]"

# Column width tweaks.
$ws.Columns.Item(2).ColumnWidth = 10.571428571428571
$ws.Columns.Item(10).ColumnWidth = 12.857142857142858

# Move the selection.
$ws.Range("I16").Select() | Out-Null
